$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "245.19"
Set-TextValue "E2" "-0.69%"
Set-TextValue "D3" "28.91"
Set-TextValue "E3" "-1.78%"
Set-TextValue "D4" "5.263"
Set-TextValue "E4" "1.87%"
Set-TextValue "D5" "0.05713"
Set-TextValue "E5" "-0.02%"
Set-TextValue "D6" "6.619"
Set-TextValue "E6" "0.37%"
Set-TextValue "D7" "3.179"
Set-TextValue "E7" "3.35%"
Set-TextValue "D8" "0.8531"
Set-TextValue "E8" "-0.53%"
Set-TextValue "D9" "0.8594"
Set-TextValue "E9" "-1.62%"
Set-TextValue "D10" "0.1368"
Set-TextValue "E10" "0.14%"
Set-TextValue "D11" "0.07043"
Set-TextValue "E11" "-0.62%"
Set-TextValue "D12" "0.03158"
Set-TextValue "E12" "10.10%"
Set-TextValue "D13" "0.09290"
Set-TextValue "D14" "0.001521"
Set-TextValue "E14" "0.34%"
Set-TextValue "D15" "0.0005968"
Set-TextValue "E15" "-94.21%"
Set-TextValue "E16" "-4.31%"
Set-TextValue "D17" "3.491"
Set-TextValue "E17" "0.17%"
Set-TextValue "D18" "2.175"
Set-TextValue "E18" "-4.58%"
Set-TextValue "D19" "0.3167"
Set-TextValue "E19" "-0.16%"
Set-TextValue "D20" "0.03324"
Set-TextValue "E20" "0.30%"
Set-TextValue "E21" "-1.79%"
Set-TextValue "D22" "3.497"
Set-TextValue "E22" "0.68%"
Set-TextValue "D23" "0.04108"
Set-TextValue "E23" "-1.46%"
Set-TextValue "E24" "-0.05%"
Set-TextValue "D25" "0.001222"
Set-TextValue "E25" "0.24%"
Set-TextValue "D26" "0.004146"
Set-TextValue "E26" "-17.77%"
Set-TextValue "E27" "-0.79%"
Set-TextValue "D28" "0.0001448"
Set-TextValue "E28" "-25.28%"
Set-TextValue "D40" "0.03767"
Set-TextValue "E40" "0.30%"
Set-TextValue "D41" "0.1066"
Set-TextValue "E41" "-0.62%"
Set-TextValue "D42" "0.003697"
Set-TextValue "E42" "-35.90%"
Set-TextValue "D43" "0.002448"
Set-TextValue "E43" "16.63%"
Set-TextValue "D44" "0.009329"
Set-TextValue "E44" "-8.96%"
Set-TextValue "D45" "0.00005310"
Set-TextValue "E45" "2.75%"
Set-TextValue "E46" "0.01%"
Set-TextValue "D47" "0.07498"
Set-TextValue "E47" "5.65%"
Set-TextValue "D48" "0.002444"
Set-TextValue "E48" "-4.94%"
Set-TextValue "D49" "0.00002099"
Set-TextValue "E49" "0.01%"
Set-TextValue "D50" "0.0001999"
Set-TextValue "E50" "0.01%"
